# Apply the updated Metrics values, let dependent "today" sheet formulas
# recalc automatically, then restore the original selections recorded in
# the diff (Metrics!C24 and today!H9).

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$today   = $wb.Worksheets.Item("today")

$metrics.Range("B2").Value  = 134656.94
$metrics.Range("B3").Value  = 117963.23999999999
$metrics.Range("B4").Value  = 41866.86
$metrics.Range("B5").Value  = 5576
$metrics.Range("B6").Value  = 4930902.6900000004
$metrics.Range("B7").Value  = 4160039.9200000009
$metrics.Range("B8").Value  = 1448826.69
$metrics.Range("B9").Value  = 191783
$metrics.Range("B10").Value = 33396283.680000007
$metrics.Range("B11").Value = 31435315.079999998
$metrics.Range("B12").Value = 11730548.730000002
$metrics.Range("B13").Value = 1289413

$excel.Calculate()

$metrics.Activate()
$metrics.Range("C24").Select()

$today.Activate()
$today.Range("H9").Select()
